$wb = $excel.ActiveWorkbook

# --- "Version History" sheet: add new v1.2 row ---
$wsHist = $wb.Worksheets.Item("Version History")
$wsHist.Range("A4").Value = "v1.2"
$wsHist.Range("B4").Value = "Ahmed Abuzaid"
$wsHist.Range("C4").Value = "verify the status after modification done on the test cases"
$wsHist.Range("D4").Formula = "=TODAY()"

# --- Reviews sheet: mark reviewer verification as Closed ---
$wsReviews = $wb.Worksheets.Item(" LH_TC_IDCONSTRAINS_REVIEWS")
$wsReviews.Range("J2:J5").Value = "Closed"
